# fix #683: replace axis "geo" by "country" in the tutorial_IO examples
# workbook (larray/tests/data/examples.xlsx). The "geo" axis header lives
# in cell A1 of every sheet (pop, births, deaths, pop_births_deaths,
# pop_missing_axis_name, pop_missing_values, pop_narrow_format).
$wb = $excel.ActiveWorkbook

# xlWhole = 1 -> only replace cells whose *entire* content is "geo",
# leaving any other string (e.g. "gender") untouched.
$xlWhole = 1

foreach ($ws in $wb.Worksheets) {
    $null = $ws.Cells.Replace("geo", "country", $xlWhole)
}
